# Update projection outputs and fix termination compensation proration logic
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New computed values per row (row -> column letter -> value)
$data = @{
    2 = @{
        B = 102; C = 102; D = 86
        E = 0.8431372549019608; F = 0.8431372549019608
        G = 0.09877977698439171; H = 0.0832849100064479
        I = 461132.027355649; J = 167952.0138788245; K = 0
        L = 167952.0138788245; M = 629084.0412344737
        N = 10240780.2488; O = 9833039.318699999
        P = 0.01640031421419328; Q = 0.01708037651791155
    }
    3 = @{
        B = 106; C = 106; D = 91
        E = 0.8584905660377359; F = 0.8584905660377359
        G = 0.09496116222475419; H = 0.08152326190993048
        I = 486228.6741479071; J = 177901.9172809835; K = 0
        L = 177901.9172809835; M = 664130.5914288907
        N = 10797211.559364; O = 10389838.401361
        P = 0.01647665383815658; Q = 0.01712268376163383
    }
    4 = @{
        B = 109; C = 109; D = 93
        E = 0.8532110091743119; F = 0.8532110091743119
        G = 0.09300393342188501; H = 0.07935197989206702
        I = 512838.5062540149; J = 184120.4607402484; K = 0
        L = 184120.4607402484; M = 696958.9669942633
        N = 11361855.35814492; O = 10953511.00540183
        P = 0.01620514035221006; Q = 0.01680926423038673
    }
    5 = @{
        B = 110; C = 109; D = 93
        E = 0.8532110091743119; F = 0.8454545454545455
        G = 0.09289780249835385; H = 0.0785408693849719
        I = 523961.5408676272; J = 187108.0142754518; K = 0
        L = 187108.0142754518; M = 711069.5551430788
        N = 11517660.85818927; O = 11106966.17486389
        P = 0.01624531374722798; Q = 0.01684600559051803
    }
    6 = @{
        B = 112; C = 112; D = 95
        E = 0.8482142857142857; F = 0.8482142857142857
        G = 0.09190722738373355; H = 0.07795702322727401
        I = 548708.848341326; J = 196420.4853275503; K = 0
        L = 196420.4853275503; M = 745129.3336688762
        N = 12112094.92793495; O = 11697629.4041098
        P = 0.01621688787086141; Q = 0.01679147787487101
    }
}

foreach ($row in $data.Keys) {
    $cols = $data[$row]
    foreach ($col in $cols.Keys) {
        $ws.Range("$col$row").Value = $cols[$col]
    }
}
